# Update "想去人数" (F) values on the "展览" and "全部类型" sheets:
#   Row 2 (丽水·动漫游戏展):        F2  432 -> 434
#   Row 3 (丽水·CCAC动漫游戏嘉年华): F3  7   -> 10

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 434
    $ws.Range("F3").Value = 10
}
